$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo in Kiss Béla's address (row 3, column E): "Ünagy" -> "Nagy"
$ws.Range("E3").Value = "Győr Nagy utca 6."

# Row 4 - Horváth Réka
$ws.Range("A4").Value = "Karcsi"
$ws.Range("B4").Value = "Horváth Réka"
$ws.Range("C4").Value = "36 30 942 3614"
$ws.Range("D4").Value = "horvath.reka@gmail.com"
$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:horvath.reka@gmail.com")
$ws.Range("D4").Style = "Hivatkozás"
$ws.Range("E4").Value = "Győr Jereváni út 23."

# Row 5 - Kovács Ernő
$ws.Range("A5").Value = "Lüszi"
$ws.Range("B5").Value = "Kovács Ernő"
$ws.Range("C5").Value = "30 70 542 7823"
$ws.Range("D5").Value = "kovacs.erno@gmail.com"
$ws.Hyperlinks.Add($ws.Range("D5"), "mailto:kovacs.erno@gmail.com")
$ws.Range("D5").Style = "Hivatkozás"
$ws.Range("E5").Value = "Győr Tihanyi Árpád utca 10."

# Row 6 - Balogh Sára
$ws.Range("A6").Value = "Artemisz"
$ws.Range("B6").Value = "Balogh Sára"
$ws.Range("C6").Value = "30 40 589 2561"
$ws.Range("D6").Value = "balogh.sara@gmail.com"
$ws.Hyperlinks.Add($ws.Range("D6"), "mailto:balogh.sara@gmail.com")
$ws.Range("D6").Style = "Hivatkozás"
$ws.Range("E6").Value = "Győr Mester utca 9."

$ws.Range("E6").Select() | Out-Null
